# "remove the % in the rows"
# Column N (header "%norm_c") is formatted as a percentage (0.058 -> "5.8%").
# Convert it to a plain number so it displays "5.8" instead of "5.8%":
#  - multiply the stored values by 100
#  - change the number format from percentage to a plain "0.00" number format

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CDF1")

# Column N ("%norm_c") holds header in row 1 and data in row 2..lastRow.
$lastRow = $ws.UsedRange.Rows.Count
$col = 14

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $cell.Value2 = [Math]::Round($cell.Value2 * 100, 10)
}

# Apply a plain numeric format (no percent sign) to the header and data cells.
$headerAndData = $ws.Range($ws.Cells.Item(1, $col), $ws.Cells.Item($lastRow, $col))
$headerAndData.NumberFormat = "0.00"
$ws.Columns.Item($col).NumberFormat = "0.00"

# Reflect the final selection recorded in the sheet view.
$ws.Range("N79").Select()
